$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "Elena Bruschetti"
$ws.Range("B56").Value = "Nicola Togni | RSA United"
$ws.Range("C56").Value = "Stefano  Galvagni | Clitoriders"
$ws.Range("D56").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("E56").Value = "ANDREA ASTE | Pinguini Trentini"
$ws.Range("F56").Value = "Davide  Bazzano  | iMontagna"
